$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1456
$ws.Range("I58").Value = 99.90000000000001
$ws.Range("J58").Value = 15017
$ws.Range("K58").Value = 299.7
$ws.Range("L58").Value = 45051
$ws.Range("M58").Value = -149.7
$ws.Range("N58").Value = -45351
$ws.Range("H64").Value = 3039.8
$ws.Range("I64").Value = 2799.6
$ws.Range("K64").Value = 2799.6
$ws.Range("M64").Value = -2551.6
$ws.Range("H67").Value = 3039.8
$ws.Range("I67").Value = 2799.6
$ws.Range("K67").Value = 2799.6
$ws.Range("M67").Value = -1941.6
$ws.Range("H74").Value = 4558.375
$ws.Range("I74").Value = 2931.125
$ws.Range("J74").Value = 5372
$ws.Range("K74").Value = 2931.125
$ws.Range("L74").Value = 5372
$ws.Range("M74").Value = -1995.125
$ws.Range("N74").Value = -7244
$ws.Range("H76").Value = 15159062
$ws.Range("I76").Value = 9795.333000000001
$ws.Range("K76").Value = 9795.333000000001
$ws.Range("M76").Value = -9480.333000000001
$ws.Range("H77").Value = 4558.375
$ws.Range("I77").Value = 2931.125
$ws.Range("J77").Value = 5372
$ws.Range("K77").Value = 14655.625
$ws.Range("L77").Value = 26860
$ws.Range("M77").Value = -9975.625
$ws.Range("N77").Value = -36220
$ws.Range("H79").Value = 15159062
$ws.Range("I79").Value = 9795.333000000001
$ws.Range("K79").Value = 9795.333000000001
$ws.Range("M79").Value = -8703.333000000001
$ws.Range("H113").Value = 4225
$ws.Range("I113").Value = 4880
$ws.Range("J113").Value = 3757.1428
$ws.Range("K113").Value = 4880
$ws.Range("L113").Value = 3757.1428
$ws.Range("M113").Value = -1626
$ws.Range("N113").Value = -10265.1428
$ws.Range("H121").Value = 794.7143
$ws.Range("I121").Value = 334.75
$ws.Range("J121").Value = 902.94116
$ws.Range("K121").Value = 1004.25
$ws.Range("L121").Value = 2708.82348
$ws.Range("M121").Value = 742.75
$ws.Range("N121").Value = -6202.82348
$ws.Range("H127").Value = 1271.5333
$ws.Range("J127").Value = 1505.1
$ws.Range("L127").Value = 4515.299999999999
$ws.Range("N127").Value = -14435.3
$ws.Range("H137").Value = 2196.1667
$ws.Range("J137").Value = 2608.8333
$ws.Range("L137").Value = 7826.499899999999
$ws.Range("N137").Value = -12926.4999
$ws.Range("H138").Value = 2996.34
$ws.Range("I138").Value = 1084.6052
$ws.Range("J138").Value = 4168.0483
$ws.Range("K138").Value = 3253.8156
$ws.Range("L138").Value = 12504.1449
$ws.Range("M138").Value = 1886.1844
$ws.Range("N138").Value = -22784.1449

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3299.5454
$ws.Range("I63").Value = 2521.6667
$ws.Range("J63").Value = 6800
$ws.Range("K63").Value = 2521.6667
$ws.Range("L63").Value = 6800
$ws.Range("M63").Value = -1835.6667
$ws.Range("N63").Value = -8172
$ws.Range("H66").Value = 3299.5454
$ws.Range("I66").Value = 2521.6667
$ws.Range("J66").Value = 6800
$ws.Range("K66").Value = 12608.3335
$ws.Range("L66").Value = 34000
$ws.Range("M66").Value = -9176.333500000001
$ws.Range("N66").Value = -40864
$ws.Range("H88").Value = 2445.3635
$ws.Range("I88").Value = 2416.5557
$ws.Range("J88").Value = 2465.3076
$ws.Range("K88").Value = 2416.5557
$ws.Range("L88").Value = 2465.3076
$ws.Range("M88").Value = -2010.5557
$ws.Range("N88").Value = -3277.3076
$ws.Range("H91").Value = 2445.3635
$ws.Range("I91").Value = 2416.5557
$ws.Range("J91").Value = 2465.3076
$ws.Range("K91").Value = 2416.5557
$ws.Range("L91").Value = 2465.3076
$ws.Range("M91").Value = -1012.5557
$ws.Range("N91").Value = -5273.3076

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 26317988
$ws.Range("I105").Value = 2227.577
$ws.Range("J105").Value = 83335470
$ws.Range("K105").Value = 2227.577
$ws.Range("L105").Value = 83335470
$ws.Range("M105").Value = -480.5770000000002
$ws.Range("N105").Value = -83338964
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").Value = ""
$ws.Range("H141").Value = 41488.43
$ws.Range("I141").Value = 41488.43
$ws.Range("K141").Value = 41488.43
$ws.Range("M141").Value = -36308.43

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1999
$ws.Range("I31").Value = 1999
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1999
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1704
$ws.Range("N31").Value = ""
$ws.Range("H34").Value = 1999
$ws.Range("I34").Value = 1999
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1999
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1797
$ws.Range("N34").Value = ""
$ws.Range("H62").Value = 3487.76
$ws.Range("J62").Value = 4090
$ws.Range("L62").Value = 4090
$ws.Range("N62").Value = -5338
$ws.Range("H65").Value = 3487.76
$ws.Range("J65").Value = 4090
$ws.Range("L65").Value = 20450
$ws.Range("N65").Value = -26690
$ws.Range("H132").Value = 7755508
$ws.Range("I132").Value = 977.52
$ws.Range("K132").Value = 2932.56
$ws.Range("M132").Value = -402.5599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 182.71428
$ws.Range("I6").Value = 95.333336
$ws.Range("J6").Value = 340
$ws.Range("K6").Value = 286.000008
$ws.Range("L6").Value = 1020
$ws.Range("M6").Value = -173.000008
$ws.Range("N6").Value = -1246
$ws.Range("H40").Value = 171.43478
$ws.Range("I40").Value = 167.6923
$ws.Range("J40").Value = 176.3
$ws.Range("K40").Value = 670.7692
$ws.Range("L40").Value = 705.2
$ws.Range("M40").Value = -601.7692
$ws.Range("N40").Value = -843.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = ""
$ws.Range("H70").Value = 4178.7144
$ws.Range("I70").Value = 4109.75
$ws.Range("K70").Value = 4109.75
$ws.Range("M70").Value = -3839.75
$ws.Range("H73").Value = 4178.7144
$ws.Range("I73").Value = 4109.75
$ws.Range("K73").Value = 4109.75
$ws.Range("M73").Value = -3173.75
$ws.Range("H80").Value = 4003151.2
$ws.Range("I80").Value = 4399.3
$ws.Range("K80").Value = 4399.3
$ws.Range("M80").Value = -3401.3
$ws.Range("H83").Value = 4003151.2
$ws.Range("I83").Value = 4399.3
$ws.Range("K83").Value = 21996.5
$ws.Range("M83").Value = -17004.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4223075.5
$ws.Range("I22").Value = 12658227
$ws.Range("J22").Value = 5500
$ws.Range("K22").Value = 12658227
$ws.Range("L22").Value = 5500
$ws.Range("M22").Value = -12657932
$ws.Range("N22").Value = -6090
$ws.Range("H27").Value = 4223075.5
$ws.Range("I27").Value = 12658227
$ws.Range("J27").Value = 5500
$ws.Range("K27").Value = 12658227
$ws.Range("L27").Value = 5500
$ws.Range("M27").Value = -12658120
$ws.Range("N27").Value = -5714
$ws.Range("H50").Value = 6000
$ws.Range("J50").Value = 6000
$ws.Range("L50").Value = 6000
$ws.Range("N50").Value = -7274
$ws.Range("H61").Value = 1729.6086
$ws.Range("I61").Value = 1388.3889
$ws.Range("J61").Value = 2958
$ws.Range("K61").Value = 1388.3889
$ws.Range("L61").Value = 2958
$ws.Range("M61").Value = -1186.3889
$ws.Range("N61").Value = -3362
$ws.Range("H113").Value = 1729.6086
$ws.Range("I113").Value = 1388.3889
$ws.Range("J113").Value = 2958
$ws.Range("K113").Value = 1388.3889
$ws.Range("L113").Value = 2958
$ws.Range("M113").Value = 781.6111000000001
$ws.Range("N113").Value = -7298
$ws.Range("H122").Value = 16141.667
$ws.Range("J122").Value = 4132.143
$ws.Range("L122").Value = 12396.429
$ws.Range("N122").Value = -17296.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 8000
$ws.Range("J42").Value = 8000
$ws.Range("L42").Value = 8000
$ws.Range("N42").Value = -8756
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").Value = ""
$ws.Range("H87").Value = 29900
$ws.Range("J87").Value = 29900
$ws.Range("L87").Value = 29900
$ws.Range("N87").Value = -32396
$ws.Range("H90").Value = 29900
$ws.Range("J90").Value = 29900
$ws.Range("L90").Value = 89700
$ws.Range("N90").Value = -102180
$ws.Range("H132").Value = 16686885
$ws.Range("I132").Value = 40042240
$ws.Range("J132").Value = 4488.8
$ws.Range("K132").Value = 120126720
$ws.Range("L132").Value = 13466.4
$ws.Range("M132").Value = -120124190
$ws.Range("N132").Value = -18526.4
